$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: Participants tab ---
$qParticipants = @'
MATCH (p:participant)-->(s:study)
OPTIONAL MATCH (samp:sample)-->(p)
OPTIONAL MATCH (p)<--(diag:diagnosis)
OPTIONAL MATCH (samp)<--(f:file)
OPTIONAL MATCH (f)<--(g:genomic_info)
WITH s, p, samp, f, g, diag
WHERE COALESCE(g.library_strategy, 'Not specified in data') in ['Not specified in data']
With p
OPTIONAL MATCH (p)-->(s:study)
OPTIONAL MATCH (samp:sample)-->(p)
WITH s, p, apoc.coll.sort(collect(distinct samp.sample_id)) as samp
RETURN 
coalesce(p.participant_id,'') as `Participant ID`,
coalesce(s.study_name, '') as `Study Name`,
coalesce(s.phs_accession,'') as `Accession`,
coalesce(p.gender,'') as `Gender`,
coalesce(apoc.text.join(samp, ','), '') as `Samples`
ORDER BY p.participant_id limit 100
'@

# --- Row 3: Samples tab ---
$qSamples = @'
MATCH (samp:sample)-->(p:participant)-->(s:study)
OPTIONAL MATCH (samp)<--(f:file)
OPTIONAL MATCH (f)<--(g:genomic_info)
OPTIONAL MATCH (p)<--(diag:diagnosis)
WITH s, p, samp, f, g, diag
WHERE COALESCE(g.library_strategy, "Not specified in data") in ['Not specified in data']
WITH DISTINCT s, p, samp
RETURN
    coalesce(samp.sample_id, '') as `Sample ID`,
    coalesce(p.participant_id,'') as `Participant ID`,
    coalesce(s.study_name, '') as `Study Name`,
    coalesce(s.phs_accession,'') as `Accession`,
    coalesce(samp.sample_tumor_status,'') as `Tumor`,
    coalesce(samp.sample_type,'') as `Analyte Type`
ORDER BY samp.sample_id limit 100
'@

# --- Row 4: Files tab ---
$qFiles = @'
MATCH (f:file)-->(s:study)
OPTIONAL MATCH (samp:sample)<--(f)
OPTIONAL MATCH (samp)-->(p:participant)
OPTIONAL MATCH (f)<--(g:genomic_info)
OPTIONAL MATCH (p)<--(diag:diagnosis)
WITH s, p, samp, f, g, diag
WHERE COALESCE(g.library_strategy, "Not specified in data") in ['Not specified in data']
WITH DISTINCT f, s, p, samp
RETURN
    coalesce(f.file_name, '') as `File Name`,
    coalesce(s.study_name,'') as `Study Name`,
    coalesce(s.phs_accession,'') as `Accession`,
    coalesce(p.participant_id, 'Not specified in data') as `Participant ID`,
    coalesce(samp.sample_id, 'Not specified in data') as `Sample ID`,
    coalesce(f.file_type, '') as `File Type`
ORDER BY f.file_name limit 100
'@

# --- Shared StatQuery column (C2:C4) ---
$qStat = @'
CALL{
    MATCH (p:participant)-->(s:study)
    OPTIONAL MATCH (samp:sample)-->(p)
    OPTIONAL MATCH (samp)<--(f:file)
    OPTIONAL MATCH (f)<--(g:genomic_info)
    OPTIONAL MATCH (p)<--(diag:diagnosis)
    WITH s, p, samp, f, g, diag
    WHERE COALESCE( g.library_strategy, "Not specified in data") in ['Not specified in data']
    RETURN 
        count(distinct p) AS num_participants
}
WITH num_participants
CALL {
    MATCH (samp:sample)-->(p:participant)-->(s)
    OPTIONAL MATCH (samp)<--(f:file)
    OPTIONAL MATCH (p)<--(diag:diagnosis)
    OPTIONAL MATCH (f)<--(g:genomic_info)
    OPTIONAL MATCH (p)<--(diag:diagnosis)
    WITH s, p, samp, f, g, diag
    WHERE COALESCE( g.library_strategy, "Not specified in data") in ['Not specified in data']
    RETURN 
        count(distinct samp) AS num_samples
}
WITH num_participants, num_samples
CALL {
    MATCH (f:file)-->(s:study)
    OPTIONAL MATCH (f)<--(g:genomic_info)
    OPTIONAL MATCH (samp:sample)<--(f)
    OPTIONAL MATCH (p:participant)<--(samp)
    OPTIONAL MATCH (p)<--(diag:diagnosis)
    WITH s, p, samp, f, g, diag
    WHERE COALESCE( g.library_strategy, "Not specified in data") in ['Not specified in data']
    RETURN 
        count(distinct s) AS num_studies,
        count(distinct f) AS num_files
}
RETURN 
    num_studies AS Studies,
    num_participants AS Participants,
    num_samples AS Samples,
    num_files AS `Files`
'@

$neo4jFile = "TC05_CDS_Filter_LibraryStrategy-NotSpecifiedinData_Neo4jData.xlsx"
$webFile = "TC05_CDS_Filter_LibraryStrategy-NotSpecifiedinData_WebData.xlsx"

# Row 2 - ParticipantsTab (was CasesTab)
$ws.Range("A2").Value = "ParticipantsTab"
$ws.Range("B2").Value = $qParticipants
$ws.Range("C2").Value = $qStat
$ws.Range("D2").Value = $neo4jFile
$ws.Range("E2").Value = $webFile

# Row 3 - SamplesTab
$ws.Range("A3").Value = "SamplesTab"
$ws.Range("B3").Value = $qSamples
$ws.Range("C3").Value = $qStat
$ws.Range("D3").Value = $neo4jFile
$ws.Range("E3").Value = $webFile

# Row 4 - FilesTab
$ws.Range("A4").Value = "FilesTab"
$ws.Range("B4").Value = $qFiles
$ws.Range("C4").Value = $qStat
$ws.Range("D4").Value = $neo4jFile
$ws.Range("E4").Value = $webFile

# Row heights to match content growth
$ws.Rows.Item(2).RowHeight = 409.5
$ws.Rows.Item(3).RowHeight = 352.5
$ws.Rows.Item(4).RowHeight = 395.25

# Column widths (best-fit approximations)
$ws.Columns.Item(1).ColumnWidth = 20.1

# View: scroll to A3 and select C4
$ws.Range("C4").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 3
$win.ScrollColumn = 1
$ws.Range("C4").Select()

Write-Host "Edit applied"
